# The document has two headers (first-page + default) and two footers
# (first-page + default), each holding a single inline picture:
#   - Headers: BTec_Logo-Orange, currently named "image2.jpg" -> rename to "image1.jpg"
#   - Footers: PearsonLogo.png,  currently named "image1.png" -> rename to "image2.png"
#
# InlineShape objects in Word's COM model don't enumerate header/footer
# stories directly, so we go through Sections(1).Headers/Footers and grab
# the lone picture out of each story's Range.

$d = $word.ActiveDocument

# --- Headers (BTec logo): image2.jpg -> image1.jpg -------------------------

$sec = $d.Sections.Item(1)
$hdrFirst = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
$hdrFirst.Name = "image1.jpg"

$sec = $d.Sections.Item(1)
$hdrDefault = $sec.Headers.Item(1).Range.InlineShapes.Item(1)
$hdrDefault.Name = "image1.jpg"

# --- Footers (Pearson logo): image1.png -> image2.png ----------------------
# (Re-fetching the shape through its own .Range.InlineShapes avoids a stale
# handle on the footer story.)

$sec = $d.Sections.Item(1)
$ftrFirst = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
$ftrFirst = $ftrFirst.Range.InlineShapes.Item(1)
$ftrFirst.Name = "image2.png"

$sec = $d.Sections.Item(1)
$ftrDefault = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
$ftrDefault = $ftrDefault.Range.InlineShapes.Item(1)
$ftrDefault.Name = "image2.png"
